$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$ws.Cells.Item(12, 1).Value = 123456789
$ws.Cells.Item(12, 2).Value = "Real Programmers Count 0123456789 From Zero"

$ws.Cells.Item(13, 1).Value = 123456789
$ws.Cells.Item(13, 2).Value = "Real Programmers Count 0123456789 From Zero"

$ws.Cells.Item(14, 1).Value = 123456789
$ws.Cells.Item(14, 2).Value = "Real Programmers Count 0123456789 From Zero"
